$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text before writing, so values like "517.94" or
# "57.865.44" are stored as literal text (matching the source data)
# instead of being auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.865.44"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "2.453.76"
$ws.Range("E3").Value = "  -2.79%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "517.94"
$ws.Range("E5").Value = "  -3.60%  "

$ws.Range("D6").Value = "131.57"
$ws.Range("E6").Value = "  -3.43%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("D9").Value = "2.458.35"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").Value = "0.0982"
$ws.Range("E10").Value = "  -3.00%  "

$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").Value = "2.884.47"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "57.787.68"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").Value = "22.17"
$ws.Range("E16").Value = "  -3.71%  "

$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -2.77%  "

$ws.Range("D18").Value = "2.452.51"
$ws.Range("E18").Value = "  -2.16%  "

$ws.Range("D19").Value = "10.65"
$ws.Range("E19").Value = "  -3.91%  "

$ws.Range("D20").Value = "4.15"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "318.04"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  -3.87%  "

$ws.Range("D24").Value = "64.22"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -2.98%  "

$ws.Range("D28").Value = "7.32"
$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").Value = "0.0₃0735"
$ws.Range("E29").Value = "  -4.27%  "

$ws.Range("D30").Value = "165.79"
$ws.Range("E30").Value = "  -3.22%  "

$ws.Range("D31").Value = "1.68"
$ws.Range("E31").Value = "  -4.40%  "

$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  -6.60%  "

$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "18.02"
$ws.Range("E36").Value = "  -1.90%  "

$ws.Range("E37").Value = "  -7.05%  "

$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  -3.18%  "

$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  -4.52%  "

$ws.Range("D40").Value = "0.784"
$ws.Range("E40").Value = "  -3.36%  "

$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  -4.76%  "

$ws.Range("D42").Value = "271.19"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").Value = "5.01"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("D44").Value = "0.588"
$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("D45").Value = "124.36"
$ws.Range("E45").Value = "  -4.47%  "

$ws.Range("D46").Value = "0.0905"
$ws.Range("E46").Value = "  -1.94%  "

$ws.Range("D47").Value = "0.0485"
$ws.Range("E47").Value = "  -3.86%  "

$ws.Range("D48").Value = "0.0209"
$ws.Range("E48").Value = "  -4.58%  "

$ws.Range("D49").Value = "16.64"
$ws.Range("E49").Value = "  -4.29%  "

$ws.Range("D50").Value = "1.721.05"
$ws.Range("E50").Value = "  -1.97%  "

$ws.Range("D51").Value = "0.967"
$ws.Range("E51").Value = "  -2.25%  "

# Restore the default (unstyled) cell style on column D so the saved
# workbook does not carry a stray explicit number format.
$dRange.Style = "Normal"
